$d = $word.ActiveDocument

$d.Content.Find.Execute("2022-12-30 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2022-12-31 Saturday", 2) | Out-Null
$d.Content.Find.Execute("54+3=", $true, $false, $false, $false, $false, $true, 1, $false, "4+23=", 2) | Out-Null
$d.Content.Find.Execute("18-1=", $true, $false, $false, $false, $false, $true, 1, $false, "25+20=", 2) | Out-Null
$d.Content.Find.Execute("72-49=", $true, $false, $false, $false, $false, $true, 1, $false, "80-2=", 2) | Out-Null
$d.Content.Find.Execute("52+7=", $true, $false, $false, $false, $false, $true, 1, $false, "81-33=", 2) | Out-Null
$d.Content.Find.Execute("94-16=", $true, $false, $false, $false, $false, $true, 1, $false, "89-38=", 2) | Out-Null
$d.Content.Find.Execute("32+57=", $true, $false, $false, $false, $false, $true, 1, $false, "71+1=", 2) | Out-Null
$d.Content.Find.Execute("33+50=", $true, $false, $false, $false, $false, $true, 1, $false, "20+79=", 2) | Out-Null
$d.Content.Find.Execute("71-7=", $true, $false, $false, $false, $false, $true, 1, $false, "30+62=", 2) | Out-Null
$d.Content.Find.Execute("64+29=", $true, $false, $false, $false, $false, $true, 1, $false, "53-14=", 2) | Out-Null
$d.Content.Find.Execute("60+22=", $true, $false, $false, $false, $false, $true, 1, $false, "17+32=", 2) | Out-Null
$d.Content.Find.Execute("24+52=", $true, $false, $false, $false, $false, $true, 1, $false, "2+54=", 2) | Out-Null
$d.Content.Find.Execute("15+77=", $true, $false, $false, $false, $false, $true, 1, $false, "70-48=", 2) | Out-Null
$d.Content.Find.Execute("40+57=", $true, $false, $false, $false, $false, $true, 1, $false, "29+69=", 2) | Out-Null
$d.Content.Find.Execute("33-19=", $true, $false, $false, $false, $false, $true, 1, $false, "65+19=", 2) | Out-Null
$d.Content.Find.Execute("41-10=", $true, $false, $false, $false, $false, $true, 1, $false, "51-29=", 2) | Out-Null
$d.Content.Find.Execute("32+49=", $true, $false, $false, $false, $false, $true, 1, $false, "60+25=", 2) | Out-Null
$d.Content.Find.Execute("46-1=", $true, $false, $false, $false, $false, $true, 1, $false, "20+63=", 2) | Out-Null
$d.Content.Find.Execute("68+18=", $true, $false, $false, $false, $false, $true, 1, $false, "4-1=", 2) | Out-Null
$d.Content.Find.Execute("46-29=", $true, $false, $false, $false, $false, $true, 1, $false, "99-79=", 2) | Out-Null
$d.Content.Find.Execute("16+33=", $true, $false, $false, $false, $false, $true, 1, $false, "14+61=", 2) | Out-Null
$d.Content.Find.Execute("77+4=", $true, $false, $false, $false, $false, $true, 1, $false, "22+5=", 2) | Out-Null
$d.Content.Find.Execute("13+18=", $true, $false, $false, $false, $false, $true, 1, $false, "85-41=", 2) | Out-Null
$d.Content.Find.Execute("25-22=", $true, $false, $false, $false, $false, $true, 1, $false, "99-97=", 2) | Out-Null
$d.Content.Find.Execute("15+59=", $true, $false, $false, $false, $false, $true, 1, $false, "62+27=", 2) | Out-Null
$d.Content.Find.Execute("32+43=", $true, $false, $false, $false, $false, $true, 1, $false, "30+22=", 2) | Out-Null
$d.Content.Find.Execute("7+87=", $true, $false, $false, $false, $false, $true, 1, $false, "71-1=", 2) | Out-Null
$d.Content.Find.Execute("68-18=", $true, $false, $false, $false, $false, $true, 1, $false, "91+3=", 2) | Out-Null
$d.Content.Find.Execute("37+9=", $true, $false, $false, $false, $false, $true, 1, $false, "93-64=", 2) | Out-Null
$d.Content.Find.Execute("76-17=", $true, $false, $false, $false, $false, $true, 1, $false, "36+54=", 2) | Out-Null
$d.Content.Find.Execute("27-0=", $true, $false, $false, $false, $false, $true, 1, $false, "85-75=", 2) | Out-Null
$d.Content.Find.Execute("26-10=", $true, $false, $false, $false, $false, $true, 1, $false, "9+3=", 2) | Out-Null
$d.Content.Find.Execute("94-72=", $true, $false, $false, $false, $false, $true, 1, $false, "56+11=", 2) | Out-Null
$d.Content.Find.Execute("92+7=", $true, $false, $false, $false, $false, $true, 1, $false, "3+31=", 2) | Out-Null
$d.Content.Find.Execute("39-16=", $true, $false, $false, $false, $false, $true, 1, $false, "40+10=", 2) | Out-Null
$d.Content.Find.Execute("3+3=", $true, $false, $false, $false, $false, $true, 1, $false, "30+64=", 2) | Out-Null
$d.Content.Find.Execute("30-4=", $true, $false, $false, $false, $false, $true, 1, $false, "62-13=", 2) | Out-Null
$d.Content.Find.Execute("65+18=", $true, $false, $false, $false, $false, $true, 1, $false, "23-5=", 2) | Out-Null
$d.Content.Find.Execute("90-83=", $true, $false, $false, $false, $false, $true, 1, $false, "23-8=", 2) | Out-Null
$d.Content.Find.Execute("61+18=", $true, $false, $false, $false, $false, $true, 1, $false, "74-72=", 2) | Out-Null
$d.Content.Find.Execute("9+34=", $true, $false, $false, $false, $false, $true, 1, $false, "81-16=", 2) | Out-Null
$d.Content.Find.Execute("70+21=", $true, $false, $false, $false, $false, $true, 1, $false, "99-20=", 2) | Out-Null
$d.Content.Find.Execute("1+89=", $true, $false, $false, $false, $false, $true, 1, $false, "84-73=", 2) | Out-Null
$d.Content.Find.Execute("95-54=", $true, $false, $false, $false, $false, $true, 1, $false, "86-74=", 2) | Out-Null
$d.Content.Find.Execute("31-25=", $true, $false, $false, $false, $false, $true, 1, $false, "69-9=", 2) | Out-Null
$d.Content.Find.Execute("32-18=", $true, $false, $false, $false, $false, $true, 1, $false, "42-34=", 2) | Out-Null
$d.Content.Find.Execute("21+55=", $true, $false, $false, $false, $false, $true, 1, $false, "90-5=", 2) | Out-Null
$d.Content.Find.Execute("97-41=", $true, $false, $false, $false, $false, $true, 1, $false, "93-29=", 2) | Out-Null
$d.Content.Find.Execute("25+15=", $true, $false, $false, $false, $false, $true, 1, $false, "66+11=", 2) | Out-Null
$d.Content.Find.Execute("31+64=", $true, $false, $false, $false, $false, $true, 1, $false, "72-72=", 2) | Out-Null
$d.Content.Find.Execute("72-57=", $true, $false, $false, $false, $false, $true, 1, $false, "84-9=", 2) | Out-Null
$d.Content.Find.Execute("18+73=", $true, $false, $false, $false, $false, $true, 1, $false, "2+55=", 2) | Out-Null
$d.Content.Find.Execute("79-18=", $true, $false, $false, $false, $false, $true, 1, $false, "27-2=", 2) | Out-Null
$d.Content.Find.Execute("99-9=", $true, $false, $false, $false, $false, $true, 1, $false, "37+11=", 2) | Out-Null
$d.Content.Find.Execute("43-31=", $true, $false, $false, $false, $false, $true, 1, $false, "25+22=", 2) | Out-Null
$d.Content.Find.Execute("4+72=", $true, $false, $false, $false, $false, $true, 1, $false, "52+39=", 2) | Out-Null
$d.Content.Find.Execute("50-8=", $true, $false, $false, $false, $false, $true, 1, $false, "76-70=", 2) | Out-Null
$d.Content.Find.Execute("79-40=", $true, $false, $false, $false, $false, $true, 1, $false, "49-43=", 2) | Out-Null
$d.Content.Find.Execute("54+36=", $true, $false, $false, $false, $false, $true, 1, $false, "21-18=", 2) | Out-Null
$d.Content.Find.Execute("28+4=", $true, $false, $false, $false, $false, $true, 1, $false, "20+74=", 2) | Out-Null
$d.Content.Find.Execute("34+53=", $true, $false, $false, $false, $false, $true, 1, $false, "83-57=", 2) | Out-Null
$d.Content.Find.Execute("52+36=", $true, $false, $false, $false, $false, $true, 1, $false, "14-13=", 2) | Out-Null
$d.Content.Find.Execute("33+1=", $true, $false, $false, $false, $false, $true, 1, $false, "43-28=", 2) | Out-Null
$d.Content.Find.Execute("40+46=", $true, $false, $false, $false, $false, $true, 1, $false, "80+1=", 2) | Out-Null
$d.Content.Find.Execute("50-12=", $true, $false, $false, $false, $false, $true, 1, $false, "17+66=", 2) | Out-Null
$d.Content.Find.Execute("13+23=", $true, $false, $false, $false, $false, $true, 1, $false, "46-16=", 2) | Out-Null
$d.Content.Find.Execute("45+5=", $true, $false, $false, $false, $false, $true, 1, $false, "17+0=", 2) | Out-Null
$d.Content.Find.Execute("21+10=", $true, $false, $false, $false, $false, $true, 1, $false, "89-49=", 2) | Out-Null
$d.Content.Find.Execute("50+19=", $true, $false, $false, $false, $false, $true, 1, $false, "55-46=", 2) | Out-Null
$d.Content.Find.Execute("3+9=", $true, $false, $false, $false, $false, $true, 1, $false, "50+39=", 2) | Out-Null
$d.Content.Find.Execute("21+27=", $true, $false, $false, $false, $false, $true, 1, $false, "60-40=", 2) | Out-Null
$d.Content.Find.Execute("56+29=", $true, $false, $false, $false, $false, $true, 1, $false, "41-18=", 2) | Out-Null
$d.Content.Find.Execute("70-8=", $true, $false, $false, $false, $false, $true, 1, $false, "74-7=", 2) | Out-Null
$d.Content.Find.Execute("79+18=", $true, $false, $false, $false, $false, $true, 1, $false, "57-31=", 2) | Out-Null
$d.Content.Find.Execute("40+13=", $true, $false, $false, $false, $false, $true, 1, $false, "80-55=", 2) | Out-Null
$d.Content.Find.Execute("33-31=", $true, $false, $false, $false, $false, $true, 1, $false, "36+44=", 2) | Out-Null
$d.Content.Find.Execute("96-15=", $true, $false, $false, $false, $false, $true, 1, $false, "15+79=", 2) | Out-Null
$d.Content.Find.Execute("58-27=", $true, $false, $false, $false, $false, $true, 1, $false, "12+58=", 2) | Out-Null
$d.Content.Find.Execute("77+17=", $true, $false, $false, $false, $false, $true, 1, $false, "99-1=", 2) | Out-Null
$d.Content.Find.Execute("74+7=", $true, $false, $false, $false, $false, $true, 1, $false, "50+22=", 2) | Out-Null
$d.Content.Find.Execute("24+2=", $true, $false, $false, $false, $false, $true, 1, $false, "69-12=", 2) | Out-Null
$d.Content.Find.Execute("74-51=", $true, $false, $false, $false, $false, $true, 1, $false, "37-24=", 2) | Out-Null
$d.Content.Find.Execute("34+26=", $true, $false, $false, $false, $false, $true, 1, $false, "95-84=", 2) | Out-Null
$d.Content.Find.Execute("34-12=", $true, $false, $false, $false, $false, $true, 1, $false, "87-61=", 2) | Out-Null
$d.Content.Find.Execute("84-76=", $true, $false, $false, $false, $false, $true, 1, $false, "63-22=", 2) | Out-Null
$d.Content.Find.Execute("8+0=", $true, $false, $false, $false, $false, $true, 1, $false, "86-35=", 2) | Out-Null
$d.Content.Find.Execute("63-44=", $true, $false, $false, $false, $false, $true, 1, $false, "77-71=", 2) | Out-Null
$d.Content.Find.Execute("10+28=", $true, $false, $false, $false, $false, $true, 1, $false, "96-58=", 2) | Out-Null
$d.Content.Find.Execute("93-53=", $true, $false, $false, $false, $false, $true, 1, $false, "93-9=", 2) | Out-Null
$d.Content.Find.Execute("46+39=", $true, $false, $false, $false, $false, $true, 1, $false, "22+67=", 2) | Out-Null
$d.Content.Find.Execute("83-63=", $true, $false, $false, $false, $false, $true, 1, $false, "59-12=", 2) | Out-Null
$d.Content.Find.Execute("88-83=", $true, $false, $false, $false, $false, $true, 1, $false, "44-17=", 2) | Out-Null
$d.Content.Find.Execute("81+14=", $true, $false, $false, $false, $false, $true, 1, $false, "47+34=", 2) | Out-Null
$d.Content.Find.Execute("82-25=", $true, $false, $false, $false, $false, $true, 1, $false, "22+14=", 2) | Out-Null
$d.Content.Find.Execute("74-49=", $true, $false, $false, $false, $false, $true, 1, $false, "58-7=", 2) | Out-Null
$d.Content.Find.Execute("11-5=", $true, $false, $false, $false, $false, $true, 1, $false, "66-64=", 2) | Out-Null
$d.Content.Find.Execute("54-10=", $true, $false, $false, $false, $false, $true, 1, $false, "19+4=", 2) | Out-Null
$d.Content.Find.Execute("9+77=", $true, $false, $false, $false, $false, $true, 1, $false, "61+17=", 2) | Out-Null
$d.Content.Find.Execute("45+34=", $true, $false, $false, $false, $false, $true, 1, $false, "28+42=", 2) | Out-Null
$d.Content.Find.Execute("4+78=", $true, $false, $false, $false, $false, $true, 1, $false, "6+25=", 2) | Out-Null
$d.Content.Find.Execute("18+51=", $true, $false, $false, $false, $false, $true, 1, $false, "64-51=", 2) | Out-Null

Write-Host "Replacements complete"
